# Apply cryptos list update (price/volume refresh) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.128.29'
$ws.Cells.Item(2, 5).Value = '  -0.65%  '

$ws.Cells.Item(3, 4).Value = '2.076.70'
$ws.Cells.Item(3, 5).Value = '  -0.95%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '253.48'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.03%  '

$ws.Cells.Item(6, 5).Value = '  +1.55%  '

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '59.17'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +7.74%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 5).Value = '  +4.62%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '61.50'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -0.63%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0802'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +7.60%  '

$ws.Cells.Item(12, 5).Value = '  +2.49%  '

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.30'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +6.96%  '

$ws.Cells.Item(14, 5).Value = '  -0.88%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.820'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -2.48%  '

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.54'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +6.74%  '

$ws.Cells.Item(17, 4).Value = '2.078.28'
$ws.Cells.Item(17, 5).Value = '  -0.99%  '

$ws.Cells.Item(18, 4).Value = '37.137.93'
$ws.Cells.Item(18, 5).Value = '  -0.52%  '

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.75'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +6.48%  '

$ws.Cells.Item(20, 5).Value = '  +2.58%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0929'
$ws.Cells.Item(21, 5).Value = '  +8.85%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.48'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +4.90%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '239.51'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.76%  '

$ws.Cells.Item(24, 5).Value = '  -0.08%  '

$ws.Cells.Item(25, 5).Value = '  -2.77%  '

$ws.Cells.Item(26, 5).Value = '  +13.74%  '

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '170.29'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.23%  '

$ws.Cells.Item(28, 5).Value = '  +0.75%  '

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.38'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -1.68%  '

$ws.Cells.Item(31, 5).Value = '  +6.36%  '

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.79'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +5.92%  '

$ws.Cells.Item(33, 5).Value = '  +3.06%  '

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.50'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +8.56%  '

$ws.Cells.Item(35, 5).Value = '  +0.38%  '

$ws.Cells.Item(36, 5).Value = '  +0.06%  '

$ws.Cells.Item(37, 5).Value = '  +2.19%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.117'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +26.27%  '

$ws.Cells.Item(39, 5).Value = '  -4.67%  '

$ws.Cells.Item(40, 5).Value = '  +2.01%  '

$ws.Cells.Item(41, 5).Value = '  +0.20%  '

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.89'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -3.12%  '

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.17'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.19%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '99.18'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.20%  '

$ws.Cells.Item(45, 2).Value = 'HuobiToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.84'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +1.45%  '

$ws.Cells.Item(46, 2).Value = 'FTXToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.31'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +6.05%  '

$ws.Cells.Item(47, 5).Value = '  +13.91%  '

$ws.Cells.Item(48, 5).Value = '  +7.67%  '

$ws.Cells.Item(49, 4).Value = '1.306.12'
$ws.Cells.Item(49, 5).Value = '  -1.22%  '

$ws.Cells.Item(50, 5).Value = '  -0.06%  '

$ws.Cells.Item(51, 5).Value = '  -0.79%  '
